$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (B, C, F) to match target layout
$ws.Columns.Item(2).ColumnWidth = 45.83333333333333
$ws.Columns.Item(3).ColumnWidth = 50.83333333333333
$ws.Columns.Item(6).ColumnWidth = 54.83333333333333

# Update timetable cell contents (columns B-F, rows 2-13)
$ws.Range("B2").Value = '{}'
$ws.Range("C2").Value = '{}'
$ws.Range("D2").Value = '{}'
$ws.Range("E2").Value = '{}'
$ws.Range("F2").Value = '{}'
$ws.Range("B3").Value = '{}'
$ws.Range("C3").Value = '{}'
$ws.Range("D3").Value = '{}'
$ws.Range("E3").Value = '{}'
$ws.Range("F3").Value = '{}'
$ws.Range("B4").Value = '{}'
$ws.Range("C4").Value = '{}'
$ws.Range("D4").Value = '{}'
$ws.Range("E4").Value = '{}'
$ws.Range("F4").Value = '{0: sala nr 8 | Jan Nowak | Język polski}'
$ws.Range("B5").Value = '{}'
$ws.Range("C5").Value = '{}'
$ws.Range("D5").Value = '{}'
$ws.Range("E5").Value = '{0: sala nr 7 | Katarzyna Mazur | Fizyka}'
$ws.Range("F5").Value = '{0: sala nr 2 | Piotr Wójcik | Biologia}'
$ws.Range("B6").Value = '{}'
$ws.Range("C6").Value = '{}'
$ws.Range("D6").Value = '{}'
$ws.Range("E6").Value = '{0: sala nr 1 | Lena Kowalska | Język angielski}'
$ws.Range("F6").Value = '{0: sala nr 11 | Dominik Kaczor | Informatyka}'
$ws.Range("B7").Value = '{}'
$ws.Range("C7").Value = '{}'
$ws.Range("D7").Value = '{}'
$ws.Range("E7").Value = '{0: sala nr 7 | Zofia Wiśniewska | Wychowanie fizyczne}'
$ws.Range("F7").Value = '{0: sala nr 4 | Karolina Kamińska | Chemia}'
$ws.Range("B8").Value = '{}'
$ws.Range("C8").Value = '{}'
$ws.Range("D8").Value = '{}'
$ws.Range("E8").Value = '{0: sala nr 5 | Paweł Lewandowski | Matematyka}'
$ws.Range("F8").Value = '{0: sala nr 10 | Dominik Kaczor | Informatyka}'
$ws.Range("B9").Value = '{}'
$ws.Range("C9").Value = '{0: sala nr 1 | Paweł Lewandowski | Matematyka}'
$ws.Range("D9").Value = '{}'
$ws.Range("E9").Value = '{}'
$ws.Range("F9").Value = '{0: sala nr 10 | Katarzyna Mazur | Fizyka}'
$ws.Range("B10").Value = '{0: sala nr 11 | Dominik Kaczor | Informatyka}'
$ws.Range("C10").Value = '{0: sala nr 9 | Piotr Wójcik | Biologia}'
$ws.Range("D10").Value = '{0: sala nr 4 | Zofia Wiśniewska | Wychowanie fizyczne}'
$ws.Range("E10").Value = '{0: sala nr 7 | Katarzyna Mazur | Fizyka}'
$ws.Range("F10").Value = '{0: sala nr 3 | Karolina Kamińska | Chemia}'
$ws.Range("B11").Value = '{0: sala nr 2 | Natalia Szymańska | Geografia}'
$ws.Range("C11").Value = '{0: sala nr 5 | Jan Nowak | Język polski}'
$ws.Range("D11").Value = '{0: sala nr 4 | Paweł Lewandowski | Matematyka}'
$ws.Range("E11").Value = '{0: sala nr 4 | Jan Nowak | Język polski}'
$ws.Range("F11").Value = '{}'
$ws.Range("B12").Value = '{0: sala nr 3 | Dominik Kaczor | Informatyka}'
$ws.Range("C12").Value = '{0: sala nr 10 | Paweł Lewandowski | Matematyka}'
$ws.Range("D12").Value = '{0: sala nr 6 | Lena Kowalska | Język angielski}'
$ws.Range("E12").Value = '{0: sala nr 8 | Natalia Szymańska | Geografia}'
$ws.Range("F12").Value = '{0: sala nr 8 | Mateusz Kowalski | Język niemiecki}'
$ws.Range("B13").Value = '{0: sala nr 5 | Dominik Kaczor | Informatyka}'
$ws.Range("C13").Value = '{0: sala nr 1 | Mateusz Kowalski | Język niemiecki}'
$ws.Range("D13").Value = '{0: sala nr 1 | Lena Kowalska | Język angielski}'
$ws.Range("E13").Value = '{0: sala nr 3 | Paweł Lewandowski | Matematyka}'
$ws.Range("F13").Value = '{0: sala nr 8 | Zofia Wiśniewska | Wychowanie fizyczne}'
